# Update "想去人数" (want-to-go headcount) figures for the 南宁-漫展信息 workbook.
# These edits mirror the upstream commit "Update gh-pages to output generated at 456a3b4",
# which bumped attendee counts on the 展览 (Exhibitions) and 演出 (Performances) sheets,
# and correspondingly on the combined 全部类型 (All types) sheet.

$wb = $excel.ActiveWorkbook

# --- 展览 (Exhibitions) sheet ---
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 6295
$wsExhibit.Range("F3").Value = 23
$wsExhibit.Range("F4").Value = 183
$wsExhibit.Range("F5").Value = 1007
$wsExhibit.Range("F6").Value = 103

# --- 演出 (Performances) sheet ---
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F2").Value = 6

# --- 全部类型 (All types) sheet -- mirrors both sheets above ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 6295
$wsAll.Range("F3").Value = 23
$wsAll.Range("F4").Value = 183
$wsAll.Range("F5").Value = 1007
$wsAll.Range("F6").Value = 103
$wsAll.Range("F7").Value = 6
